$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "equip compose" row (Id=4, row 7 on the sheet) is being removed as part
# of optimising the equip compose system. Deleting the entire row shifts the
# rows below it up by one, which also shrinks the table range (handled
# automatically since the sheet data lives inside a structured table).
$ws.Rows.Item(7).Delete()

# Re-select row 7 (now holding what used to be row 8's data) to match the
# workbook's recorded selection state after the edit.
$ws.Range("A7:XFD7").Select()
